$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The title cell A1 ("APURAÇÃO DO CUMPRIMENTO DO LIMITE LEGAL") is removed.
# Clearing it drops the now-unused shared string and shifts nothing else --
# the remaining row-1 header cells (B1:F1) stay put.
$ws.Range("A1").ClearContents() | Out-Null

# With the long title gone, column A's widest content becomes the row
# labels in A3:A12, so Excel auto-sizes the column to fit them.
$ws.Columns.Item(1).ColumnWidth = 120.8

# The saved view now shows cell B4 selected.
$ws.Range("B4").Select() | Out-Null
